$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 212-213 (pushes the old 212..325 block down to 214..327).
$ws.Rows("212:213").Insert()

# New "Primera" record (row 212) for the week of 2022-02-04.
$ws.Cells.Item(212, 1).Value = 8
$ws.Cells.Item(212, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(212, 3).Value = "Coquimbo"
$ws.Cells.Item(212, 4).Value = 44596
$ws.Cells.Item(212, 5).Value = 4
$ws.Cells.Item(212, 6).Value = 100112009
$ws.Cells.Item(212, 7).Value = "Acelga"
$ws.Cells.Item(212, 8).Value = "Sin especificar"
$ws.Cells.Item(212, 9).Value = "Primera"
$ws.Cells.Item(212, 10).Value = 2200
$ws.Cells.Item(212, 11).Value = 450
$ws.Cells.Item(212, 12).Value = 500
$ws.Cells.Item(212, 13).Value = 475
$ws.Cells.Item(212, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(212, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(212, 16).Value = 238
$ws.Cells.Item(212, 17).Value = 2
$ws.Cells.Item(212, 18).Value = "Hortaliza"

# New "Segunda" record (row 213) for the same week.
$ws.Cells.Item(213, 1).Value = 8
$ws.Cells.Item(213, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(213, 3).Value = "Coquimbo"
$ws.Cells.Item(213, 4).Value = 44596
$ws.Cells.Item(213, 5).Value = 4
$ws.Cells.Item(213, 6).Value = 100112009
$ws.Cells.Item(213, 7).Value = "Acelga"
$ws.Cells.Item(213, 8).Value = "Sin especificar"
$ws.Cells.Item(213, 9).Value = "Segunda"
$ws.Cells.Item(213, 10).Value = 1320
$ws.Cells.Item(213, 11).Value = 350
$ws.Cells.Item(213, 12).Value = 400
$ws.Cells.Item(213, 13).Value = 375
$ws.Cells.Item(213, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(213, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(213, 16).Value = 188
$ws.Cells.Item(213, 17).Value = 2
$ws.Cells.Item(213, 18).Value = "Hortaliza"
